$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.158.23"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.42"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7391"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.09"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07221"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.94"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08354"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7613"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.923.90"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.466"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.14"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.205.11"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.163"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "250.53"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.66"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007899"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.184.57"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.971"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1580"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.324"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.75"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.78"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.065"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.485"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.592"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.536"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.211"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05374"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.259"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7825"
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.005"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.764"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4576"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.099.77"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.090"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.83"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8732"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.53"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.867"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.611"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.082.83"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.625"
$ws.Range("E51").Value = "  -1.06%  "
